$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / date row ---
$ws.Range('A1').Value = 'C6 Brawler Levels'
$ws.Range('A2').Value = '2022-12-29'

# --- Row 2 hidden summary cells (top player stats), same as row 4 ---
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = '#JVQP8LLY'
$ws.Cells.Item(2, 4).Value = 37212
$ws.Cells.Item(2, 5).Value = 10
$ws.Cells.Item(2, 6).Value = 'Diamond I'
$ws.Cells.Item(2, 7).Value = 26
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 34
$ws.Cells.Item(2, 10).Value = 'BEA, BELLE, BIBI, BO, BROCK, BUSTER, BUZZ, BYRON, CARL, CHESTER, COLT, CROW, ... '

# --- Data rows 4-33 ---
$rows = @(
  @{ Row=4; A='SuperGod'; B=1; C='#JVQP8LLY'; D=37212; E=10; F='Diamond I'; G=26; H=4; I=34; J='BEA, BELLE, BIBI, BO, BROCK, BUSTER, BUZZ, BYRON, CARL, CHESTER, COLT, CROW, ... ' },
  @{ Row=5; A='Ys'; B=1; C='#C00LQCUP'; D=30543; E=6; F='Silver III'; G=38; H=11; I=14; J='BARLEY, BO, CROW, DARRYL, EMZ, EVE, JESSIE, NANI, NITA, PENNY, POCO, SANDY, ... ' },
  @{ Row=6; A='cuteangel'; B=1; C='#JJV99PV9'; D=37894; E=13; F='Mythic I'; G=24; H=6; I=33; J='AMBER, ASH, BEA, BELLE, BONNIE, BULL, BUSTER, CARL, CHESTER, COLETTE, CROW, EVE, ... ' },
  @{ Row=7; A='Lava Blaze'; B=2; C='#2QPRUQJP2'; D=25032; E=10; F='Diamond I'; G=21; H=3; I=8; J='DARRYL, EL PRIMO, GALE, JESSIE, LEON, SANDY, SPIKE, SURGE' },
  @{ Row=8; A='MR.'; B=2; C='#9V8GGUQRR'; D=34951; E=12; F='Diamond III'; G=22; H=6; I=28; J='ASH, BELLE, BO, BONNIE, BROCK, BUSTER, BUZZ, CARL, COLETTE, CROW, DARRYL, EMZ, ... ' },
  @{ Row=9; A='Tsunami splash'; B=2; C='#82GLQ8VQJ'; D=37251; E=11; F='Diamond II'; G=27; H=17; I=15; J='ASH, BEA, BELLE, BROCK, BUSTER, CARL, CHESTER, COLETTE, CROW, FANG, MORTIS, SPIKE, ... ' },
  @{ Row=10; A='RayDaBest'; B=3; C='#200U2RLLV'; D=36795; E=13; F='Mythic I'; G=18; H=32; I=4; J='8-BIT, ASH, BIBI, TICK' },
  @{ Row=11; A='komji'; B=3; C='#8J9V9U02L'; D=40070; E=12; F='Diamond III'; G=0; H=32; I=31; J='AMBER, ASH, BARLEY, BO, BONNIE, BROCK, BULL, BUSTER, BUZZ, BYRON, CARL, CROW, ... ' },
  @{ Row=12; A='宇智波鼬'; B=3; C='#VR2Q28UG'; D=30079; E=8; F='Gold II'; G=38; H=6; I=10; J='BEA, BO, EMZ, GUS, JACKY, MORTIS, PENNY, PIPER, POCO, SQUEAK' },
  @{ Row=13; A='KMB | ɴᴀᴛᴛɪᴇ♡'; B=4; C='#2CQYQU0QQ'; D=12714; E=7; F='Gold I'; G=3; H=39; I=2; J='BUZZ, PENNY' },
  @{ Row=14; A='_mym_'; B=4; C='#82808UG9G'; D=35688; E=10; F='Diamond I'; G=30; H=16; I=18; J='BEA, BO, BUSTER, BYRON, CARL, CHESTER, CROW, GENE, GRIFF, MAX, PAM, POCO, ... ' },
  @{ Row=15; A='nik haikal'; B=4; C='#P00C0RR8'; D=30070; E=12; F='Diamond III'; G=2; H=27; I=34; J='8-BIT, BEA, BELLE, BO, BONNIE, BROCK, BUSTER, BUZZ, CARL, CHESTER, COLT, CROW, ... ' },
  @{ Row=16; A='axnsan'; B=5; C='#2VY2PC0PL'; D=33410; E=9; F='Gold III'; G=17; H=15; I=32; J='8-BIT, AMBER, ASH, BELLE, BO, BONNIE, BROCK, BYRON, CHESTER, COLETTE, CROW, EDGAR, ... ' },
  @{ Row=17; A='lolzorsish'; B=5; C='#289GU8LR8'; D=30291; E=13; F='Mythic I'; G=33; H=14; I=15; J='BO, BROCK, CARL, COLT, DARRYL, FRANK, GALE, GENE, LEON, NITA, PAM, POCO, ... ' },
  @{ Row=18; A='xardas'; B=5; C='#2P88VGRL0'; D=33767; E=13; F='Mythic I'; G=22; H=0; I=41; J='BEA, BELLE, BO, BONNIE, BROCK, BUSTER, BUZZ, CARL, CHESTER, COLETTE, CROW, DARRYL, ... ' },
  @{ Row=19; A='LA | FLASH'; B=6; C='#2YCQJ00Y'; D=32051; E=13; F='Mythic I'; G=37; H=8; I=9; J='COLETTE, COLT, CROW, FANG, MAX, MORTIS, STU, SURGE, TARA' },
  @{ Row=20; A='Saurav'; B=6; C='#U2Q9L2QU'; D=34855; E=14; F='Mythic II'; G=22; H=24; I=10; J='BEA, CARL, COLT, EMZ, GENE, JACKY, MAX, MORTIS, SURGE, TARA' },
  @{ Row=21; A='Snoopy>.<|をゆひせぬ'; B=6; C='#C9RCCU8J'; D=35134; E=11; F='Diamond II'; G=5; H=43; I=12; J='8-BIT, BUSTER, CHESTER, DYNAMIKE, EDGAR, EMZ, GRAY, MEG, PENNY, POCO, STU' },
  @{ Row=22; A='Blaze'; B=7; C='#C0R8YQC'; D=32136; E=10; F='Diamond I'; G=27; H=18; I=14; J='BELLE, BUSTER, BYRON, CHESTER, CROW, GRIFF, GUS, JANET, OTIS, PIPER, POCO, SQUEAK, ... ' },
  @{ Row=23; A='DOOM'; B=7; C='#CPJC0QUV'; D=36353; E=11; F='Diamond II'; G=20; H=12; I=31; J='ASH, BEA, BELLE, BO, BONNIE, BUSTER, CARL, CHESTER, CROW, EMZ, EVE, GENE, ... ' },
  @{ Row=24; A='eric'; B=7; C='#80VR8V9'; D=33025; E=13; F='Mythic I'; G=26; H=18; I=13; J='ASH, BELLE, EVE, GENE, GRIFF, LOLA, LOU, NANI, POCO, SPIKE, SPROUT, SURGE, ... ' },
  @{ Row=25; A='Tribe | LHC 2'; B=8; C='#V8VRPRYQ'; D=20183; E=16; F='Legend I'; G=2; H=4; I=20; J='8-BIT, BONNIE, BROCK, CARL, CROW, DARRYL, EL PRIMO, EMZ, EVE, FANG, GROM, GUS, ... ' },
  @{ Row=26; A='joshua'; B=8; C='#C29RQJLU'; D=39987; E=8; F='Gold II'; G=0; H=40; I=23; J='BARLEY, BELLE, BIBI, BO, BYRON, COLT, CROW, DARRYL, EMZ, EVE, GRIFF, LOU, ... ' },
  @{ Row=27; A='☬ℝ𝔸𝕋𝕋𝕃𝔼ℝ☬'; B=8; C='#89GV9UG9Q'; D=40293; E=10; F='Diamond I'; G=20; H=8; I=35; J='AMBER, ASH, BEA, BELLE, BIBI, BO, BROCK, BUZZ, CARL, CHESTER, COLETTE, CROW, ... ' },
  @{ Row=28; A='IX|LIT'; B=9; C='#8V09Y2Y8'; D=32994; E=7; F='Gold I'; G=35; H=14; I=5; J='BIBI, CHESTER, MORTIS, RICO, TARA' },
  @{ Row=29; A='Mini Breeze'; B=9; C='#2CVYPV0YP'; D=9099; E=10; F='Diamond I'; G=22; H=6; I=2; J='MORTIS, RICO' },
  @{ Row=30; A='Mini|Benn🎯'; B=9; C='#9RVV02QQ'; D=7652; E=7; F='Gold I'; G=5; H=1; I=7; J='BO, DARRYL, LOLA, PAM, PENNY, RICO, SPIKE' },
  @{ Row=31; A='RICOFTW'; B=10; C='#9ULG0RR8V'; D=5554; E=16; F='Legend I'; G=0; H=0; I=7; J='BEA, COLT, GROM, JESSIE, PENNY, RICO, SHELLY' },
  @{ Row=32; A='Synderella'; B=10; C='#R80JC998'; D=32569; E=13; F='Mythic I'; G=0; H=0; I=63; J='8-BIT, AMBER, ASH, BARLEY, BEA, BELLE, BIBI, BO, BONNIE, BROCK, BULL, BUSTER, ... ' },
  @{ Row=33; A='shauntws'; B=10; C='#GJCYYV0P'; D=31761; E=11; F='Diamond II'; G=36; H=17; I=11; J='ASH, BYRON, CARL, EMZ, EVE, GENE, GRIFF, SANDY, SPIKE, SURGE' }
)

foreach ($rowData in $rows) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H
    $ws.Cells.Item($r, 9).Value = $rowData.I
    $ws.Cells.Item($r, 10).Value = $rowData.J
}

# --- New row 33: apply the same formatting pattern as the rest of the data rows ---
for ($c = 1; $c -le 9; $c++) {
    $ws.Cells.Item(33, $c).HorizontalAlignment = -4108
}
$ws.Cells.Item(33, 10).HorizontalAlignment = -4131

# --- Extend conditional formatting range to include the new row ---
$ws.Range('A2:H32').FormatConditions.Delete()
$fc = $ws.Range('A2:H33').FormatConditions.Add(2, 3, 'MOD($B2,2)=0')
